$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$a1 = $ws.Range("A1")

# The title/instructions cell (A1, merged A1:K1) contains two rich-text runs:
#  run 1 (bold)     - "DO NOT DELETE THIS ROW! ... ignored.\n"
#  run 2 (regular)  - "\nNAME: Full Name\n...YEAR OF BAPTISM: Optional...\n...false."
# We need to fix the typo/inconsistency: "YEAR OF BAPTISM" -> "YEAR OF BIRTH"
# (the actual column header in E2 already says "YEAR OF BIRTH").

$fullText = $a1.Text
$boldMarker = "ignored.`n"
$boldEnd = $fullText.IndexOf($boldMarker) + $boldMarker.Length

$target = "YEAR OF BAPTISM"
$idx = $fullText.IndexOf($target)

if ($idx -ge 0) {
    $chars = $a1.Characters($idx + 1, $target.Length)
    $chars.Text = "YEAR OF BIRTH"
}

# Replacing text via Characters collapses rich-text run formatting, so
# restore the original two-run structure (bold header, regular body).
$newFullText = $a1.Text

$run1 = $a1.Characters(1, $boldEnd)
$run1.Font.Name = "Calibri"
$run1.Font.Size = 12
$run1.Font.Bold = $true

$run2 = $a1.Characters($boldEnd + 1, $newFullText.Length - $boldEnd)
$run2.Font.Name = "Calibri"
$run2.Font.Size = 12
$run2.Font.Bold = $false

# The user also re-selected the merged title cell before saving.
[void]$ws.Range("A1:K1").Select()

# Touch page setup (paper size / orientation) as in the authored workbook.
$ps = $ws.PageSetup
$ps.PaperSize = 9
$ps.Orientation = 1
